$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "Poisson error 1/50 tolerence"
$ws.Range("C11").Value = "SOR 1.1 used"

$ws.Range("K4").Select()
